$d = $word.ActiveDocument
$tbl = $d.Tables.Item(1)

$cell = $tbl.Cell(1, 1)
$cell.Range.Find.Execute("571÷8=71, 3", $true, $false, $false, $false, $false, $true, 1, $false, "464÷8=58, 0", 1) | Out-Null
$cell = $tbl.Cell(1, 2)
$cell.Range.Find.Execute("345÷7=49, 2", $true, $false, $false, $false, $false, $true, 1, $false, "708÷4=177, 0", 1) | Out-Null
$cell = $tbl.Cell(1, 3)
$cell.Range.Find.Execute("570÷3=190, 0", $true, $false, $false, $false, $false, $true, 1, $false, "552÷7=78, 6", 1) | Out-Null
$cell = $tbl.Cell(1, 4)
$cell.Range.Find.Execute("107÷5=21, 2", $true, $false, $false, $false, $false, $true, 1, $false, "782÷9=86, 8", 1) | Out-Null
$cell = $tbl.Cell(1, 5)
$cell.Range.Find.Execute("593÷4=148, 1", $true, $false, $false, $false, $false, $true, 1, $false, "453÷6=75, 3", 1) | Out-Null
$cell = $tbl.Cell(5, 1)
$cell.Range.Find.Execute("664÷7=94, 6", $true, $false, $false, $false, $false, $true, 1, $false, "182÷5=36, 2", 1) | Out-Null
$cell = $tbl.Cell(5, 2)
$cell.Range.Find.Execute("954÷8=119, 2", $true, $false, $false, $false, $false, $true, 1, $false, "264÷5=52, 4", 1) | Out-Null
$cell = $tbl.Cell(5, 3)
$cell.Range.Find.Execute("862÷2=431, 0", $true, $false, $false, $false, $false, $true, 1, $false, "677÷5=135, 2", 1) | Out-Null
$cell = $tbl.Cell(5, 5)
$cell.Range.Find.Execute("448÷6=74, 4", $true, $false, $false, $false, $false, $true, 1, $false, "190÷2=95, 0", 1) | Out-Null
$cell = $tbl.Cell(9, 1)
$cell.Range.Find.Execute("202÷2=101, 0", $true, $false, $false, $false, $false, $true, 1, $false, "466÷7=66, 4", 1) | Out-Null
$cell = $tbl.Cell(9, 2)
$cell.Range.Find.Execute("624÷5=124, 4", $true, $false, $false, $false, $false, $true, 1, $false, "628÷4=157, 0", 1) | Out-Null
$cell = $tbl.Cell(9, 3)
$cell.Range.Find.Execute("165÷5=33, 0", $true, $false, $false, $false, $false, $true, 1, $false, "544÷3=181, 1", 1) | Out-Null
$cell = $tbl.Cell(9, 4)
$cell.Range.Find.Execute("861÷4=215, 1", $true, $false, $false, $false, $false, $true, 1, $false, "842÷3=280, 2", 1) | Out-Null
$cell = $tbl.Cell(9, 5)
$cell.Range.Find.Execute("262÷5=52, 2", $true, $false, $false, $false, $false, $true, 1, $false, "762÷3=254, 0", 1) | Out-Null
$cell = $tbl.Cell(13, 1)
$cell.Range.Find.Execute("607÷2=303, 1", $true, $false, $false, $false, $false, $true, 1, $false, "284÷9=31, 5", 1) | Out-Null
$cell = $tbl.Cell(13, 2)
$cell.Range.Find.Execute("572÷3=190, 2", $true, $false, $false, $false, $false, $true, 1, $false, "336÷6=56, 0", 1) | Out-Null
$cell = $tbl.Cell(13, 3)
$cell.Range.Find.Execute("560÷6=93, 2", $true, $false, $false, $false, $false, $true, 1, $false, "159÷6=26, 3", 1) | Out-Null
$cell = $tbl.Cell(13, 4)
$cell.Range.Find.Execute("337÷5=67, 2", $true, $false, $false, $false, $false, $true, 1, $false, "845÷3=281, 2", 1) | Out-Null
$cell = $tbl.Cell(13, 5)
$cell.Range.Find.Execute("617÷4=154, 1", $true, $false, $false, $false, $false, $true, 1, $false, "127÷7=18, 1", 1) | Out-Null
$cell = $tbl.Cell(17, 1)
$cell.Range.Find.Execute("829÷5=165, 4", $true, $false, $false, $false, $false, $true, 1, $false, "149÷7=21, 2", 1) | Out-Null
$cell = $tbl.Cell(17, 2)
$cell.Range.Find.Execute("924÷6=154, 0", $true, $false, $false, $false, $false, $true, 1, $false, "535÷4=133, 3", 1) | Out-Null
$cell = $tbl.Cell(17, 3)
$cell.Range.Find.Execute("862÷2=431, 0", $true, $false, $false, $false, $false, $true, 1, $false, "840÷8=105, 0", 1) | Out-Null
$cell = $tbl.Cell(17, 4)
$cell.Range.Find.Execute("234÷4=58, 2", $true, $false, $false, $false, $false, $true, 1, $false, "875÷8=109, 3", 1) | Out-Null
$cell = $tbl.Cell(17, 5)
$cell.Range.Find.Execute("681÷2=340, 1", $true, $false, $false, $false, $false, $true, 1, $false, "158÷7=22, 4", 1) | Out-Null
